# "added 4wk low sales check"
#
# Forecast Comparison sheet: the Inventory Coverage (H) figures are no
# longer reliable once recent sales are too thin, so those cells are
# blanked out; Stockout Risk / Reorder Urgency get recomputed off the new
# low-sales check (High/Urgent -> Low/Normal for the affected weeks), and
# the Seasonality Index (L) is recalculated.
#
# Summary sheet: with the low-sales guard in place the forward forecast
# totals collapse to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Inventory Coverage (H2:H17) is no longer computable -> blank it out.
$ws.Range("H2:H17").ClearContents()

# Stockout Risk (I) drops from High to Low for weeks 6-17.
$ws.Range("I6").Value  = "Low"
$ws.Range("I7").Value  = "Low"
$ws.Range("I8").Value  = "Low"
$ws.Range("I9").Value  = "Low"
$ws.Range("I10").Value = "Low"
$ws.Range("I11").Value = "Low"
$ws.Range("I12").Value = "Low"
$ws.Range("I13").Value = "Low"
$ws.Range("I14").Value = "Low"
$ws.Range("I15").Value = "Low"
$ws.Range("I16").Value = "Low"
$ws.Range("I17").Value = "Low"

# Reorder Urgency (J) drops from Urgent to Normal for weeks 5-17.
$ws.Range("J5").Value  = "Normal"
$ws.Range("J6").Value  = "Normal"
$ws.Range("J7").Value  = "Normal"
$ws.Range("J8").Value  = "Normal"
$ws.Range("J9").Value  = "Normal"
$ws.Range("J10").Value = "Normal"
$ws.Range("J11").Value = "Normal"
$ws.Range("J12").Value = "Normal"
$ws.Range("J13").Value = "Normal"
$ws.Range("J14").Value = "Normal"
$ws.Range("J15").Value = "Normal"
$ws.Range("J16").Value = "Normal"
$ws.Range("J17").Value = "Normal"

# Seasonality Index (L) recalculated for every week.
$ws.Range("L2").Value  = 0.86
$ws.Range("L3").Value  = 1.06
$ws.Range("L4").Value  = 1
$ws.Range("L5").Value  = 1.06
$ws.Range("L6").Value  = 0.93
$ws.Range("L7").Value  = 0.89
$ws.Range("L8").Value  = 1.17
$ws.Range("L9").Value  = 0.83
$ws.Range("L10").Value = 1.18
$ws.Range("L11").Value = 1.04
$ws.Range("L12").Value = 1
$ws.Range("L13").Value = 0.94
$ws.Range("L14").Value = 0.89
$ws.Range("L15").Value = 1.1
$ws.Range("L16").Value = 1.14
$ws.Range("L17").Value = 1.1

# Summary sheet: forward forecast totals are now all 0. These cells hold
# text-formatted numbers (same convention as the rest of the sheet, e.g.
# "0 units" on row 8), so prefix with an apostrophe to force text and keep
# Excel from re-typing them as numeric.
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").Value  = "'0"
$summary.Range("B10").Value = "'0"
$summary.Range("B11").Value = "'0"
$summary.Range("B12").Value = "'0"
